# "Generate Report for Handback"
#
# The CI handback run completed successfully for both locales (zh-cn and
# de-de): the per-locale Status moves from "Ready for handoff" to
# "Handed back: in sync with en-US", the Error Detail column (which held
# the stale "version mismatch" warning) is cleared now that the handback
# is in sync, and the Latest Handback DateTime is refreshed to the new
# run's timestamp. The Status/Error columns are also widened slightly so
# the longer message text is readable.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: Status is mirrored per-locale in columns E (zh-cn) and
# F (de-de) on row 2.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText
$overview.Range("E1").ColumnWidth = 29.166666666666668
$overview.Range("F1").ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------
# zh-cn sheet: Status (C2), Latest Handback DateTime (K2) and Error
# Detail (P2) all reflect the freshly completed handback.
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $statusText
$zhcn.Range("K2").Value = "2016-08-22 06:48:26"
$zhcn.Range("P2").Value = ""
$zhcn.Range("C1").ColumnWidth = 29.166666666666668
$zhcn.Range("P1").ColumnWidth = 12.833333333333334

# ---------------------------------------------------------------------
# de-de sheet: same shape of update as zh-cn, with its own handback
# timestamp.
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $statusText
$dede.Range("K2").Value = "2016-08-22 06:48:33"
$dede.Range("P2").Value = ""
$dede.Range("C1").ColumnWidth = 29.166666666666668
$dede.Range("P1").ColumnWidth = 12.833333333333334
